$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 and A3 with the new combined text values
$ws.Range("A2").Value = "('Angel', ['Token Creature — Angel', 'Flying', '4/4'])"
$ws.Range("A3").Value = "('Elemental', ['Token Creature — Elemental', '3/1'])"

# Delete old rows 4 through 8 which are no longer needed, shifting cells up
$ws.Range("A4:A8").EntireRow.Delete() | Out-Null
